$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Matières enseignés"

$ws.Columns.Item(3).ColumnWidth = 26.666666666666668
$ws.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws.Columns.Item(5).ColumnWidth = 30.833333333333332

$ws.Range("E6").Select()
